$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1472.5385
$ws.Range("J17").Value = 1472.5385
$ws.Range("L17").Value = 4417.6155
$ws.Range("N17").Value = -4753.6155
$ws.Range("H33").Value = 363.36365
$ws.Range("I33").Value = 363.36365
$ws.Range("K33").Value = 363.36365
$ws.Range("M33").Value = -134.36365
$ws.Range("H41").Value = 833.8570999999999
$ws.Range("I41").Value = 868
$ws.Range("K41").Value = 868
$ws.Range("M41").Value = -428
$ws.Range("H62").Value = 4065
$ws.Range("I62").Value = 4078.8
$ws.Range("K62").Value = 4078.8
$ws.Range("M62").Value = -3454.8
$ws.Range("H65").Value = 4065
$ws.Range("I65").Value = 4078.8
$ws.Range("K65").Value = 20394
$ws.Range("M65").Value = -17274
$ws.Range("H129").Value = 568.125
$ws.Range("I129").Value = 568.125
$ws.Range("K129").Value = 1704.375
$ws.Range("M129").Value = 3295.625
$ws.Range("H131").Value = 1506.5
$ws.Range("I131").Value = 1007.8
$ws.Range("K131").Value = 3023.4
$ws.Range("M131").Value = 2016.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 46650
$ws.Range("J96").Value = 46650
$ws.Range("L96").Value = 46650
$ws.Range("N96").Value = -52142
$ws.Range("H132").Value = 6997.5713
$ws.Range("I132").Value = 4984
$ws.Range("K132").Value = 14952
$ws.Range("M132").Value = -12422
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H134").Value = 4999.8335
$ws.Range("I134").Value = 4999.8335
$ws.Range("K134").Value = 14999.5005
$ws.Range("M134").Value = -12464.5005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 100001
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H37").Value = 24993
$ws.Range("I37").Value = 24994.5
$ws.Range("K37").Value = 24994.5
$ws.Range("M37").Value = -24887.5
$ws.Range("H50").Value = 5500
$ws.Range("I50").Value = 5500
$ws.Range("K50").Value = 5500
$ws.Range("M50").Value = -4875
$ws.Range("H51").Value = 3199.5
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 36701.332
$ws.Range("I59").Value = 36701.332
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 36701.332
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -35556.332
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H61").Value = 3199.5
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 7249.75
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 7999.6665
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 7999.6665
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -10245.6665
$ws.Range("H89").Value = 7249.75
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 7999.6665
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 39998.3325
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -51230.3325
$ws.Range("H99").Value = 2900
$ws.Range("I99").Value = 2900
$ws.Range("K99").Value = 2900
$ws.Range("M99").Value = -1402
$ws.Range("H122").Value = 2166.6667
$ws.Range("I122").Value = 2250
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6750
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4300
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 2900
$ws.Range("I126").Value = 2900
$ws.Range("K126").Value = 8700
$ws.Range("M126").Value = -6230
$ws.Range("H132").Value = 2076.2727
$ws.Range("I132").Value = 2112.9
$ws.Range("K132").Value = 6338.700000000001
$ws.Range("M132").Value = -3808.700000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 185.77777
$ws.Range("I2").Value = 69.818184
$ws.Range("J2").Value = 265.5
$ws.Range("K2").Value = 418.909104
$ws.Range("L2").Value = 1593
$ws.Range("M2").Value = -305.909104
$ws.Range("N2").Value = -1819
$ws.Range("H5").Value = 1174.5
$ws.Range("J5").Value = 1233
$ws.Range("L5").Value = 3699
$ws.Range("N5").Value = -3923
$ws.Range("H12").Value = 1180.625
$ws.Range("J12").Value = 1962.3334
$ws.Range("L12").Value = 5887.0002
$ws.Range("N12").Value = -6233.0002
$ws.Range("H26").Value = 1582.125
$ws.Range("J26").Value = 1866
$ws.Range("L26").Value = 5598
$ws.Range("N26").Value = -6174
$ws.Range("H50").Value = 2627.5
$ws.Range("I50").Value = 255
$ws.Range("K50").Value = 765
$ws.Range("M50").Value = -284
$ws.Range("H53").Value = 2627.5
$ws.Range("I53").Value = 255
$ws.Range("K53").Value = 765
$ws.Range("M53").Value = -284
$ws.Range("H92").Value = 298
$ws.Range("J92").Value = 246
$ws.Range("L92").Value = 738
$ws.Range("N92").Value = -3234
$ws.Range("H135").Value = 1174.5
$ws.Range("J135").Value = 1233
$ws.Range("L135").Value = 11097
$ws.Range("N135").Value = -16167
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2608
$ws.Range("I102").Value = 2608
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2608
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -986
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 3178.25
$ws.Range("I122").Value = 2584.6667
$ws.Range("J122").Value = 3771.8333
$ws.Range("K122").Value = 7754.000100000001
$ws.Range("L122").Value = 11315.4999
$ws.Range("M122").Value = -5304.000100000001
$ws.Range("N122").Value = -16215.4999
$ws.Range("H132").Value = 1969.8572
$ws.Range("I132").Value = 1958.4
$ws.Range("J132").Value = 1998.5
$ws.Range("K132").Value = 5875.200000000001
$ws.Range("L132").Value = 5995.5
$ws.Range("M132").Value = -3345.200000000001
$ws.Range("N132").Value = -11055.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3100
$ws.Range("J82").Value = 4000
$ws.Range("L82").Value = 4000
$ws.Range("N82").Value = -4722
$ws.Range("H85").Value = 3100
$ws.Range("J85").Value = 4000
$ws.Range("L85").Value = 4000
$ws.Range("N85").Value = -6496
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24995
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H45").Value = 39992.555
$ws.Range("I45").Value = 32748.25
$ws.Range("J45").Value = 45788
$ws.Range("K45").Value = 32748.25
$ws.Range("L45").Value = 45788
$ws.Range("M45").Value = -32257.25
$ws.Range("N45").Value = -46770
$ws.Range("H75").Value = 55333
$ws.Range("J75").Value = 55333
$ws.Range("L75").Value = 55333
$ws.Range("N75").Value = -57205
$ws.Range("H78").Value = 55333
$ws.Range("J78").Value = 55333
$ws.Range("L78").Value = 165999
$ws.Range("N78").Value = -175359
$ws.Range("H107").Value = 9000
$ws.Range("I107").Value = 8000
$ws.Range("J107").Value = 9333.333000000001
$ws.Range("K107").Value = 24000
$ws.Range("L107").Value = 27999.999
$ws.Range("M107").Value = -22080
$ws.Range("N107").Value = -31839.999
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("M131").ClearContents()
